# began script to block study/test orders
# Swap the paired stimulus-pairing values in column B for the six
# row-pairs that were re-shuffled: (3,6) (20,71) (154,176) (193,198)
# (195,199) (196,200).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B3").Value = 6
$ws.Range("B6").Value = 3

$ws.Range("B20").Value = 71
$ws.Range("B71").Value = 20

$ws.Range("B154").Value = 176
$ws.Range("B176").Value = 154

$ws.Range("B193").Value = 198
$ws.Range("B195").Value = 199
$ws.Range("B196").Value = 200
$ws.Range("B198").Value = 193
$ws.Range("B199").Value = 195
$ws.Range("B200").Value = 196

# Scroll/selection moved down near the bottom of the list (row ~181,
# active cell B198) as the script-writer worked through the later rows.
$ws.Range("B198").Select()
